$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'231.40"
$ws.Range("D3").Value = "'22.71"
$ws.Range("D4").Value = "'5.272"
$ws.Range("D5").Value = "'0.05591"
$ws.Range("D6").Value = "'3.377"
$ws.Range("D7").Value = "'6.462"
$ws.Range("D8").Value = "'1.060"
$ws.Range("D9").Value = "'0.7821"
$ws.Range("D10").Value = "'0.1375"
$ws.Range("D11").Value = "'0.07335"
$ws.Range("D12").Value = "'0.03138"
$ws.Range("D13").Value = "'0.02968"
$ws.Range("D14").Value = "'0.09265"
$ws.Range("D15").Value = "'0.001657"
$ws.Range("D17").Value = "'0.04752"
$ws.Range("D18").Value = "'0.0005795"
$ws.Range("E18").Value = "17OneONEWorstin24h"
$ws.Range("D19").Value = "'0.006251"
$ws.Range("D20").Value = "'0.005240"
$ws.Range("D21").Value = "'0.001054"
$ws.Range("D23").Value = "'3.971"
$ws.Range("D26").Value = "'0.1053"
$ws.Range("D27").Value = "'0.0004995"
$ws.Range("D40").Value = "'0.04015"
$ws.Range("D41").Value = "'0.006998"
$ws.Range("D42").Value = "'0.1040"
$ws.Range("D43").Value = "'0.003219"
$ws.Range("D44").Value = "'0.009783"
$ws.Range("D45").Value = "'0.00005435"
$ws.Range("D46").Value = "'0.00000000749"
$ws.Range("D47").Value = "'0.7845"
$ws.Range("D48").Value = "'0.04216"
$ws.Range("E48").Value = "47BOLOBOLO"
$ws.Range("D49").Value = "'0.00002098"
$ws.Range("D50").Value = "'0.01009"
